$d = $word.ActiveDocument

# Phase 1: replace each original value with a unique placeholder
# to avoid collisions where a new value contains another old value as a substring.
$d.Content.Find.Execute("2023-11-23 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0000@@", 2) | Out-Null
$d.Content.Find.Execute("64+2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0001@@", 2) | Out-Null
$d.Content.Find.Execute("11+47=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0002@@", 2) | Out-Null
$d.Content.Find.Execute("22+53=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0003@@", 2) | Out-Null
$d.Content.Find.Execute("91-32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0004@@", 2) | Out-Null
$d.Content.Find.Execute("32+20=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0005@@", 2) | Out-Null
$d.Content.Find.Execute("69-11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0006@@", 2) | Out-Null
$d.Content.Find.Execute("62-61=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0007@@", 2) | Out-Null
$d.Content.Find.Execute("18+34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0008@@", 2) | Out-Null
$d.Content.Find.Execute("4+26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0009@@", 2) | Out-Null
$d.Content.Find.Execute("43+11=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0010@@", 2) | Out-Null
$d.Content.Find.Execute("55+39=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0011@@", 2) | Out-Null
$d.Content.Find.Execute("10-10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0012@@", 2) | Out-Null
$d.Content.Find.Execute("50+0=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0013@@", 2) | Out-Null
$d.Content.Find.Execute("79+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0014@@", 2) | Out-Null
$d.Content.Find.Execute("28-24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0015@@", 2) | Out-Null
$d.Content.Find.Execute("56-27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0016@@", 2) | Out-Null
$d.Content.Find.Execute("54-50=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0017@@", 2) | Out-Null
$d.Content.Find.Execute("84-24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0018@@", 2) | Out-Null
$d.Content.Find.Execute("38+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0019@@", 2) | Out-Null
$d.Content.Find.Execute("22+22=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0020@@", 2) | Out-Null
$d.Content.Find.Execute("59+23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0021@@", 2) | Out-Null
$d.Content.Find.Execute("64-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0022@@", 2) | Out-Null
$d.Content.Find.Execute("8+71=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0023@@", 2) | Out-Null
$d.Content.Find.Execute("28-15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0024@@", 2) | Out-Null
$d.Content.Find.Execute("34+21=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0025@@", 2) | Out-Null
$d.Content.Find.Execute("3+61=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0026@@", 2) | Out-Null
$d.Content.Find.Execute("14+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0027@@", 2) | Out-Null
$d.Content.Find.Execute("87-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0028@@", 2) | Out-Null
$d.Content.Find.Execute("28+40=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0029@@", 2) | Out-Null
$d.Content.Find.Execute("63+32=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0030@@", 2) | Out-Null
$d.Content.Find.Execute("88-26=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0031@@", 2) | Out-Null
$d.Content.Find.Execute("79-46=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0032@@", 2) | Out-Null
$d.Content.Find.Execute("75-75=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0033@@", 2) | Out-Null
$d.Content.Find.Execute("99-6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0034@@", 2) | Out-Null
$d.Content.Find.Execute("60-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0035@@", 2) | Out-Null
$d.Content.Find.Execute("27+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0036@@", 2) | Out-Null
$d.Content.Find.Execute("28+24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0037@@", 2) | Out-Null
$d.Content.Find.Execute("5+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0038@@", 2) | Out-Null
$d.Content.Find.Execute("45+31=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0039@@", 2) | Out-Null
$d.Content.Find.Execute("27+58=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0040@@", 2) | Out-Null
$d.Content.Find.Execute("61+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0041@@", 2) | Out-Null
$d.Content.Find.Execute("98-65=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0042@@", 2) | Out-Null
$d.Content.Find.Execute("44-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0043@@", 2) | Out-Null
$d.Content.Find.Execute("51-2=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0044@@", 2) | Out-Null
$d.Content.Find.Execute("4+34=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0045@@", 2) | Out-Null
$d.Content.Find.Execute("94+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0046@@", 2) | Out-Null
$d.Content.Find.Execute("23+73=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0047@@", 2) | Out-Null
$d.Content.Find.Execute("48-27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0048@@", 2) | Out-Null
$d.Content.Find.Execute("25+51=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0049@@", 2) | Out-Null
$d.Content.Find.Execute("83-8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0050@@", 2) | Out-Null
$d.Content.Find.Execute("94-81=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0051@@", 2) | Out-Null
$d.Content.Find.Execute("9+42=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0052@@", 2) | Out-Null
$d.Content.Find.Execute("49-38=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0053@@", 2) | Out-Null
$d.Content.Find.Execute("81-23=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0054@@", 2) | Out-Null
$d.Content.Find.Execute("9+28=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0055@@", 2) | Out-Null
$d.Content.Find.Execute("34-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0056@@", 2) | Out-Null
$d.Content.Find.Execute("25+52=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0057@@", 2) | Out-Null
$d.Content.Find.Execute("62+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0058@@", 2) | Out-Null
$d.Content.Find.Execute("43+16=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0059@@", 2) | Out-Null
$d.Content.Find.Execute("87-82=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0060@@", 2) | Out-Null
$d.Content.Find.Execute("93-79=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0061@@", 2) | Out-Null
$d.Content.Find.Execute("42+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0062@@", 2) | Out-Null
$d.Content.Find.Execute("13+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0063@@", 2) | Out-Null
$d.Content.Find.Execute("33-31=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0064@@", 2) | Out-Null
$d.Content.Find.Execute("10+36=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0065@@", 2) | Out-Null
$d.Content.Find.Execute("6+67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0066@@", 2) | Out-Null
$d.Content.Find.Execute("89-67=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0067@@", 2) | Out-Null
$d.Content.Find.Execute("36-35=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0068@@", 2) | Out-Null
$d.Content.Find.Execute("12-12=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0069@@", 2) | Out-Null
$d.Content.Find.Execute("90-72=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0070@@", 2) | Out-Null
$d.Content.Find.Execute("11+43=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0071@@", 2) | Out-Null
$d.Content.Find.Execute("82+3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0072@@", 2) | Out-Null
$d.Content.Find.Execute("63-63=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0073@@", 2) | Out-Null
$d.Content.Find.Execute("12-5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0074@@", 2) | Out-Null
$d.Content.Find.Execute("38+4=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0075@@", 2) | Out-Null
$d.Content.Find.Execute("2+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0076@@", 2) | Out-Null
$d.Content.Find.Execute("64-17=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0077@@", 2) | Out-Null
$d.Content.Find.Execute("53-13=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0078@@", 2) | Out-Null
$d.Content.Find.Execute("67+25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0079@@", 2) | Out-Null
$d.Content.Find.Execute("71+7=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0080@@", 2) | Out-Null
$d.Content.Find.Execute("2+89=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0081@@", 2) | Out-Null
$d.Content.Find.Execute("55-45=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0082@@", 2) | Out-Null
$d.Content.Find.Execute("15-6=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0083@@", 2) | Out-Null
$d.Content.Find.Execute("17+71=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0084@@", 2) | Out-Null
$d.Content.Find.Execute("84-79=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0085@@", 2) | Out-Null
$d.Content.Find.Execute("56+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0086@@", 2) | Out-Null
$d.Content.Find.Execute("27+10=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0087@@", 2) | Out-Null
$d.Content.Find.Execute("42-25=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0088@@", 2) | Out-Null
$d.Content.Find.Execute("5-3=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0089@@", 2) | Out-Null
$d.Content.Find.Execute("45-1=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0090@@", 2) | Out-Null
$d.Content.Find.Execute("37+49=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0091@@", 2) | Out-Null
$d.Content.Find.Execute("66+24=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0092@@", 2) | Out-Null
$d.Content.Find.Execute("1+88=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0093@@", 2) | Out-Null
$d.Content.Find.Execute("12+15=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0094@@", 2) | Out-Null
$d.Content.Find.Execute("89+8=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0095@@", 2) | Out-Null
$d.Content.Find.Execute("44+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0096@@", 2) | Out-Null
$d.Content.Find.Execute("65+5=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0097@@", 2) | Out-Null
$d.Content.Find.Execute("18+27=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0098@@", 2) | Out-Null
$d.Content.Find.Execute("7+76=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0099@@", 2) | Out-Null
$d.Content.Find.Execute("93-91=", $true, $false, $false, $false, $false, $true, 1, $false, "@@PH0100@@", 2) | Out-Null

# Phase 2: replace each placeholder with the final new value
$d.Content.Find.Execute("@@PH0000@@", $true, $false, $false, $false, $false, $true, 1, $false, "2023-11-24 Friday", 2) | Out-Null
$d.Content.Find.Execute("@@PH0001@@", $true, $false, $false, $false, $false, $true, 1, $false, "15+64=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0002@@", $true, $false, $false, $false, $false, $true, 1, $false, "48-5=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0003@@", $true, $false, $false, $false, $false, $true, 1, $false, "40-27=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0004@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-27=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0005@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0006@@", $true, $false, $false, $false, $false, $true, 1, $false, "15+32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0007@@", $true, $false, $false, $false, $false, $true, 1, $false, "17+79=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0008@@", $true, $false, $false, $false, $false, $true, 1, $false, "43+14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0009@@", $true, $false, $false, $false, $false, $true, 1, $false, "85+1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0010@@", $true, $false, $false, $false, $false, $true, 1, $false, "59-59=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0011@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0012@@", $true, $false, $false, $false, $false, $true, 1, $false, "53-4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0013@@", $true, $false, $false, $false, $false, $true, 1, $false, "30+7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0014@@", $true, $false, $false, $false, $false, $true, 1, $false, "21+32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0015@@", $true, $false, $false, $false, $false, $true, 1, $false, "65+9=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0016@@", $true, $false, $false, $false, $false, $true, 1, $false, "49+13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0017@@", $true, $false, $false, $false, $false, $true, 1, $false, "91-29=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0018@@", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0019@@", $true, $false, $false, $false, $false, $true, 1, $false, "6+7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0020@@", $true, $false, $false, $false, $false, $true, 1, $false, "57+21=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0021@@", $true, $false, $false, $false, $false, $true, 1, $false, "80-44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0022@@", $true, $false, $false, $false, $false, $true, 1, $false, "39-22=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0023@@", $true, $false, $false, $false, $false, $true, 1, $false, "22+38=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0024@@", $true, $false, $false, $false, $false, $true, 1, $false, "23+8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0025@@", $true, $false, $false, $false, $false, $true, 1, $false, "4+75=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0026@@", $true, $false, $false, $false, $false, $true, 1, $false, "67-22=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0027@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-73=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0028@@", $true, $false, $false, $false, $false, $true, 1, $false, "63+15=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0029@@", $true, $false, $false, $false, $false, $true, 1, $false, "71+1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0030@@", $true, $false, $false, $false, $false, $true, 1, $false, "10+13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0031@@", $true, $false, $false, $false, $false, $true, 1, $false, "78-6=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0032@@", $true, $false, $false, $false, $false, $true, 1, $false, "71-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0033@@", $true, $false, $false, $false, $false, $true, 1, $false, "7+12=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0034@@", $true, $false, $false, $false, $false, $true, 1, $false, "57-33=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0035@@", $true, $false, $false, $false, $false, $true, 1, $false, "54-8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0036@@", $true, $false, $false, $false, $false, $true, 1, $false, "64+28=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0037@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+52=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0038@@", $true, $false, $false, $false, $false, $true, 1, $false, "25+63=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0039@@", $true, $false, $false, $false, $false, $true, 1, $false, "68-3=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0040@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-29=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0041@@", $true, $false, $false, $false, $false, $true, 1, $false, "20-1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0042@@", $true, $false, $false, $false, $false, $true, 1, $false, "61+34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0043@@", $true, $false, $false, $false, $false, $true, 1, $false, "34+59=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0044@@", $true, $false, $false, $false, $false, $true, 1, $false, "70-52=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0045@@", $true, $false, $false, $false, $false, $true, 1, $false, "65-46=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0046@@", $true, $false, $false, $false, $false, $true, 1, $false, "74-34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0047@@", $true, $false, $false, $false, $false, $true, 1, $false, "69-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0048@@", $true, $false, $false, $false, $false, $true, 1, $false, "27+56=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0049@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+19=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0050@@", $true, $false, $false, $false, $false, $true, 1, $false, "39+42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0051@@", $true, $false, $false, $false, $false, $true, 1, $false, "67-10=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0052@@", $true, $false, $false, $false, $false, $true, 1, $false, "2+26=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0053@@", $true, $false, $false, $false, $false, $true, 1, $false, "58+20=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0054@@", $true, $false, $false, $false, $false, $true, 1, $false, "46+16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0055@@", $true, $false, $false, $false, $false, $true, 1, $false, "93-45=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0056@@", $true, $false, $false, $false, $false, $true, 1, $false, "88-47=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0057@@", $true, $false, $false, $false, $false, $true, 1, $false, "47+51=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0058@@", $true, $false, $false, $false, $false, $true, 1, $false, "34+55=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0059@@", $true, $false, $false, $false, $false, $true, 1, $false, "9+83=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0060@@", $true, $false, $false, $false, $false, $true, 1, $false, "21-13=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0061@@", $true, $false, $false, $false, $false, $true, 1, $false, "76-7=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0062@@", $true, $false, $false, $false, $false, $true, 1, $false, "83-34=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0063@@", $true, $false, $false, $false, $false, $true, 1, $false, "85-42=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0064@@", $true, $false, $false, $false, $false, $true, 1, $false, "49-16=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0065@@", $true, $false, $false, $false, $false, $true, 1, $false, "96-89=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0066@@", $true, $false, $false, $false, $false, $true, 1, $false, "23-10=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0067@@", $true, $false, $false, $false, $false, $true, 1, $false, "31+8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0068@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-24=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0069@@", $true, $false, $false, $false, $false, $true, 1, $false, "62-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0070@@", $true, $false, $false, $false, $false, $true, 1, $false, "36+30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0071@@", $true, $false, $false, $false, $false, $true, 1, $false, "38+25=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0072@@", $true, $false, $false, $false, $false, $true, 1, $false, "11+63=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0073@@", $true, $false, $false, $false, $false, $true, 1, $false, "40+44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0074@@", $true, $false, $false, $false, $false, $true, 1, $false, "74-31=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0075@@", $true, $false, $false, $false, $false, $true, 1, $false, "84+14=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0076@@", $true, $false, $false, $false, $false, $true, 1, $false, "95-8=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0077@@", $true, $false, $false, $false, $false, $true, 1, $false, "87-68=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0078@@", $true, $false, $false, $false, $false, $true, 1, $false, "73-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0079@@", $true, $false, $false, $false, $false, $true, 1, $false, "37+36=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0080@@", $true, $false, $false, $false, $false, $true, 1, $false, "41-17=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0081@@", $true, $false, $false, $false, $false, $true, 1, $false, "42-30=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0082@@", $true, $false, $false, $false, $false, $true, 1, $false, "36+44=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0083@@", $true, $false, $false, $false, $false, $true, 1, $false, "34-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0084@@", $true, $false, $false, $false, $false, $true, 1, $false, "76+15=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0085@@", $true, $false, $false, $false, $false, $true, 1, $false, "63-56=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0086@@", $true, $false, $false, $false, $false, $true, 1, $false, "94-37=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0087@@", $true, $false, $false, $false, $false, $true, 1, $false, "18+24=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0088@@", $true, $false, $false, $false, $false, $true, 1, $false, "47-1=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0089@@", $true, $false, $false, $false, $false, $true, 1, $false, "90-55=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0090@@", $true, $false, $false, $false, $false, $true, 1, $false, "31+54=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0091@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-32=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0092@@", $true, $false, $false, $false, $false, $true, 1, $false, "12-10=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0093@@", $true, $false, $false, $false, $false, $true, 1, $false, "42+53=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0094@@", $true, $false, $false, $false, $false, $true, 1, $false, "81-75=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0095@@", $true, $false, $false, $false, $false, $true, 1, $false, "34+43=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0096@@", $true, $false, $false, $false, $false, $true, 1, $false, "1+98=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0097@@", $true, $false, $false, $false, $false, $true, 1, $false, "41+4=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0098@@", $true, $false, $false, $false, $false, $true, 1, $false, "86-43=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0099@@", $true, $false, $false, $false, $false, $true, 1, $false, "29+2=", 2) | Out-Null
$d.Content.Find.Execute("@@PH0100@@", $true, $false, $false, $false, $false, $true, 1, $false, "33-29=", 2) | Out-Null
